# Version 2.0. Failed Scenarios Runner is added.
# Populates the "Status" column (G) with pass/skip results for each
# scenario row, and rewrites the "Price" column (C) as formatted-text
# currency values instead of numeric currency-formatted values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (C): store the currency amounts as plain text -------
# Switch the cells to a text number format first so the "$" strings we
# assign are kept verbatim instead of being re-parsed back into numbers.
$ws.Range("C2:C6").NumberFormat = "@"

$ws.Range("C6").Value = "$16.40"
$ws.Range("C5").Value = "$28.98"
$ws.Range("C4").Value = "$26.00"
$ws.Range("C3").Value = "$27.00"
$ws.Range("C2").Value = "$16.51"

# Resize column C to fit its new text content.
$ws.Columns.Item(3).AutoFit()

# --- Status column (G): results of the scenario runner ----------------
$ws.Range("G2").Value = "skipped"
$ws.Range("G3").Value = "skipped"
$ws.Range("G4").Value = "passed"
$ws.Range("G5").Value = "passed"
$ws.Range("G6").Value = "passed"

# Reflect the newly populated range as the current selection.
$ws.Range("G2:G6").Select()
